$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column T (20th column) data for year 2023
$ws.Range("T4").Value = 2023
$ws.Range("T5").Value = 40
